# Remove "Academy type and route" from the ab (academy transfer) document
# template, and spell out "AB" as "Advisory Board" in the neighbouring
# "Date of AB" label.
#
# Table 2 of the template ("Project name" / "Sponsor name" /
# "Academy type and route") loses its third row entirely - that row
# (label "Academy type and route" + placeholder "[AcademyTypeAndRoute]")
# is no longer part of the template.

$d = $word.ActiveDocument

# --- 1. Delete the "Academy type and route" row -------------------------
# It is the 3rd row of the 2nd table in the document.
$projectTable = $d.Tables.Item(2)
$projectTable.Rows.Item(3).Delete()

# After the row is removed, Word re-levels the two remaining cells in the
# first column so their explicit width matches the (unchanged) grid
# column width of 5225 dxa (was 5224 dxa while the 3rd row was present).
$firstColCell = $projectTable.Rows.Item(1).Cells.Item(1)
$firstColCell.Width = 5225 / 20.0

# --- 2. Spell out "AB" as "Advisory Board" -------------------------------
# 3rd table in the document holds the "Date of AB" / "[DateOfHtb]" row;
# expand the abbreviation in that row's label cell.
$dateTable = $d.Tables.Item(3)
$labelCell = $dateTable.Rows.Item(2).Cells.Item(1)
$labelCell.Range.Find.Execute("AB", $false, $true, $false, $false, $false, `
                               $true, 1, $false, "Advisory Board", 2)
